$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")

# Row 31, column C: remove "R58" from the resistor reference designator list.
# The cell holds rich/colored text; target only the run that contains
# ",R50,R51,R57,R58," so the surrounding colored runs are left alone.
$c31 = $ws.Range("C31")
$full31 = $c31.Value()
$old31 = ",R50,R51,R57,R58,"
$new31 = ",R50,R51,R57,"
$pos31 = $full31.IndexOf($old31)
$chars31 = $c31.Characters($pos31 + 1, $old31.Length)
$chars31.Text = $new31

# Row 36, column C: remove "R56" (trailing designator) from the plain list.
$ws.Range("C36").Value = "R11,R14,R35,R36,R37,R38,R48,R49,R55"

# Update the active selection to match the saved workbook state (C34).
$ws.Range("C34").Select()
